$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clarify the PASS/FAIL status formulas: distinguish an ERROR (bad/missing
# expected value) from a genuine FAIL, instead of lumping both into FAIL.
$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'
$ws.Range("D4:D7").Formula = '=IF(ISERROR(B4),"ERROR",IF(ISERROR(C4),"FAIL",IF(B4=C4,"PASS","FAIL")))'

# Refresh the recorded QuantLib add-in call counters (re-run of the unit
# test suite produced new call-count suffixes).
$ws.Range("B5").Value = "ex01#0000"
$ws.Range("B6").Value = "ex02#0000"
$ws.Range("B7").Value = "ex03#0000"

$wb.Save()
